$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A. This shifts:
#   old A (labels)          -> B
#   old B (RawActivations)  -> C
#   old C (PercActivations) -> D
#   old D (totalActivation) -> E
$ws.Columns.Item(1).Insert()

# Give the new B1 header the same look as the other header cells, then
# set its text to "segments".
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "segments"

# Fill the new column A (rows 2-20) with the 0-based segment index, and
# copy the formatting from column B (the old column A) down onto it.
for ($i = 2; $i -le 20; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 2
}
$ws.Range("B2:B20").Copy()
$ws.Range("A2:A20").PasteSpecial(-4122)

$excel.CutCopyMode = $false
